$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: insert the new CFAP221 row at spreadsheet row 64 ---
# (shifts old rows 64-68 -- CFAP54, GOLGA3, ITCH, NME8, STK36 -- down to rows 65-69)
$ws.Rows.Item(64).Insert()

# restore the column-A bold/bordered "index" style on the newly inserted row
# (Insert() leaves the new row with a blank style)
$ws.Range("A63").Copy()
$ws.Range("A64").PasteSpecial(-4122)

# fill in the new row (gene inserted alphabetically between AKNA and CFAP54)
$ws.Range("A64").Value = 62
$ws.Range("B64").Value = "CFAP221"
$ws.Range("C64").Value = "cilia and flagella associated protein 221"
# paste D63s text-typed "1" as VALUES ONLY so D64 stays Text (matches rest of column D)
# instead of Excel auto-converting a plain "1" literal into a Number
$ws.Range("D63").Copy()
$ws.Range("D64").PasteSpecial(-4163)
$ws.Range("E64").Value = "Ciliary Dyskinesia"

# renumber the sequential index in column A for the rows pushed down by the insert
# (Insert() moves the cells but does not touch their numeric contents)
$ws.Range("A65").Value = 63
$ws.Range("A66").Value = 64
$ws.Range("A67").Value = 65
$ws.Range("A68").Value = 66
$ws.Range("A69").Value = 67

# --- Step 2: add the new time_taken column (F) ---
$ws.Range("F1").Value = "time_taken"
# copy formatting (bold font + border) from the neighbouring header cell E1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$timeTaken = @(
    "2021-10-05 10:50:28.195390",
    "2021-10-05 10:50:28.195402",
    "2021-10-05 10:50:28.195405",
    "2021-10-05 10:50:28.195408",
    "2021-10-05 10:50:28.195411",
    "2021-10-05 10:50:28.195414",
    "2021-10-05 10:50:28.195416",
    "2021-10-05 10:50:28.195419",
    "2021-10-05 10:50:28.195421",
    "2021-10-05 10:50:28.195424",
    "2021-10-05 10:50:28.195427",
    "2021-10-05 10:50:28.195429",
    "2021-10-05 10:50:28.195432",
    "2021-10-05 10:50:28.195434",
    "2021-10-05 10:50:28.195437",
    "2021-10-05 10:50:28.195439",
    "2021-10-05 10:50:28.195442",
    "2021-10-05 10:50:28.195444",
    "2021-10-05 10:50:28.195447",
    "2021-10-05 10:50:28.195449",
    "2021-10-05 10:50:28.195452",
    "2021-10-05 10:50:28.195454",
    "2021-10-05 10:50:28.195457",
    "2021-10-05 10:50:28.195459",
    "2021-10-05 10:50:28.195462",
    "2021-10-05 10:50:28.195465",
    "2021-10-05 10:50:28.195467",
    "2021-10-05 10:50:28.195470",
    "2021-10-05 10:50:28.195472",
    "2021-10-05 10:50:28.195475",
    "2021-10-05 10:50:28.195477",
    "2021-10-05 10:50:28.195479",
    "2021-10-05 10:50:28.195482",
    "2021-10-05 10:50:28.195485",
    "2021-10-05 10:50:28.195487",
    "2021-10-05 10:50:28.195490",
    "2021-10-05 10:50:28.195492",
    "2021-10-05 10:50:28.195495",
    "2021-10-05 10:50:28.195497",
    "2021-10-05 10:50:28.195499",
    "2021-10-05 10:50:28.195502",
    "2021-10-05 10:50:28.195505",
    "2021-10-05 10:50:28.195507",
    "2021-10-05 10:50:28.195510",
    "2021-10-05 10:50:28.195512",
    "2021-10-05 10:50:28.195514",
    "2021-10-05 10:50:28.195517",
    "2021-10-05 10:50:28.195519",
    "2021-10-05 10:50:28.195522",
    "2021-10-05 10:50:28.195524",
    "2021-10-05 10:50:28.195527",
    "2021-10-05 10:50:28.195529",
    "2021-10-05 10:50:28.195532",
    "2021-10-05 10:50:28.195534",
    "2021-10-05 10:50:28.195537",
    "2021-10-05 10:50:28.195539",
    "2021-10-05 10:50:28.195542",
    "2021-10-05 10:50:28.195544",
    "2021-10-05 10:50:28.195547",
    "2021-10-05 10:50:28.195549",
    "2021-10-05 10:50:28.195552",
    "2021-10-05 10:50:28.195554",
    "2021-10-05 10:50:28.195557",
    "2021-10-05 10:50:28.195559",
    "2021-10-05 10:50:28.195563",
    "2021-10-05 10:50:28.195565",
    "2021-10-05 10:50:28.195568",
    "2021-10-05 10:50:28.195570"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $timeTaken[$i]
}

